{"js": "const replacements = [\n  [\"41\u00d732=\", \"97\u00d786=\"],\n  [\"32\u00d781=\", \"29\u00d794=\"],\n  [\"36\u00d733=\", \"31\u00d739=\"],\n  [\"91\u00d771=\", \"12\u00d763=\"],\n  [\"44\u00d795=\", \"19\u00d717=\"],\n  [\"66\u00d763=\", \"73\u00d748=\"],\n  [\"56\u00d726=\", \"16\u00d787=\"],\n  [\"91\u00d797=\", \"91\u00d738=\"],\n  [\"26\u00d734=\", \"16\u00d788=\"],\n  [\"53\u00d737=\", \"32\u00d743=\"],\n  [\"61\u00d713=\", \"47\u00d745=\"],\n  [\"79\u00d744=\", \"14\u00d775=\"],\n  [\"16\u00d785=\", \"56\u00d797=\"],\n  [\"99\u00d795=\", \"38\u00d724=\"],\n  [\"27\u00d781=\", \"77\u00d784=\"],\n  [\"72\u00d753=\", \"20\u00d735=\"],\n  [\"13\u00d789=\", \"84\u00d769=\"],\n  [\"83\u00d789=\", \"54\u00d795=\"],\n  [\"84\u00d791=\", \"23\u00d747=\"],\n  [\"78\u00d734=\", \"57\u00d785=\"],\n  [\"11\u00d749=\", \"62\u00d721=\"],\n  [\"81\u00d731=\", \"34\u00d752=\"],\n  [\"53\u00d760=\", \"63\u00d713=\"],\n  [\"17\u00d785=\", \"73\u00d731=\"],\n  [\"30\u00d765=\", \"62\u00d734=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @{\n    \"41\u00d732=\" = \"97\u00d786=\"\n    \"32\u00d781=\" = \"29\u00d794=\"\n    \"36\u00d733=\" = \"31\u00d739=\"\n    \"91\u00d771=\" = \"12\u00d763=\"\n    \"44\u00d795=\" = \"19\u00d717=\"\n    \"66\u00d763=\" = \"73\u00d748=\"\n    \"56\u00d726=\" = \"16\u00d787=\"\n    \"91\u00d797=\" = \"91\u00d738=\"\n    \"26\u00d734=\" = \"16\u00d788=\"\n    \"53\u00d737=\" = \"32\u00d743=\"\n    \"61\u00d713=\" = \"47\u00d745=\"\n    \"79\u00d744=\" = \"14\u00d775=\"\n    \"16\u00d785=\" = \"56\u00d797=\"\n    \"99\u00d795=\" = \"38\u00d724=\"\n    \"27\u00d781=\" = \"77\u00d784=\"\n    \"72\u00d753=\" = \"20\u00d735=\"\n    \"13\u00d789=\" = \"84\u00d769=\"\n    \"83\u00d789=\" = \"54\u00d795=\"\n    \"84\u00d791=\" = \"23\u00d747=\"\n    \"78\u00d734=\" = \"57\u00d785=\"\n    \"11\u00d749=\" = \"62\u00d721=\"\n    \"81\u00d731=\" = \"34\u00d752=\"\n    \"53\u00d760=\" = \"63\u00d713=\"\n    \"17\u00d785=\" = \"73\u00d731=\"\n    \"30\u00d765=\" = \"62\u00d734=\"\n}\n\nforeach ($key in $replacements.Keys) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($key, $false, $false, $false, $false, $false, $true, 1, $false, $replacements[$key], 2)\n}\n"}
